$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 4 PerManova tables each get a new "Total R2" column (K) summing G (R2) values
# for that block. Header cell gets the new shared string + centered style; the first
# data row of each block gets a SUM formula over the block's R2 column.
$blocks = @(
    @{ Header = 2;  DataStart = 3;  DataEnd = 8  },
    @{ Header = 12; DataStart = 13; DataEnd = 18 },
    @{ Header = 22; DataStart = 23; DataEnd = 28 },
    @{ Header = 32; DataStart = 33; DataEnd = 38 }
)

foreach ($b in $blocks) {
    $headerCell = "K" + $b.Header
    $ws.Range($headerCell).Value = "Total R2"
    $ws.Range($headerCell).HorizontalAlignment = -4108  # xlCenter

    $formulaCell = "K" + $b.DataStart
    $sumRange = "G" + $b.DataStart + ":G" + $b.DataEnd
    $ws.Range($formulaCell).Formula = "=SUM(" + $sumRange + ")"
}

# Update the window scroll position / selection to match the new view state.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("G28").Select()
